$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price and volume data
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '40.125.66'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.222.01'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '294.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '87.84'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.08%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '30.62'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '50.87'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0781'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('E13').Value = '  +3.86%  '
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.566.97'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.82'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.225.36'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.735'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '40.061.91'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0889'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.03%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.65'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '235.66'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('E28').Value = '  +3.47%  '
$ws.Range('E29').Value = '  +1.78%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.07'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.61'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.82'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.88%  '
$ws.Range('E34').Value = '  +0.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.04'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.32%  '
$ws.Range('E36').Value = '  -0.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.71%  '
$ws.Range('E38').Value = '  +1.73%  '
$ws.Range('E39').Value = '  +4.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '15.62'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.075.04'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.50%  '
$ws.Range('E43').Value = '  -2.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.59'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +13.30%  '
$ws.Range('E45').Value = '  +1.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.77'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.437.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('E50').Value = '  +5.19%  '
$ws.Range('E51').Value = '  +2.01%  '
